$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
# Row 18
$ws.Range("H18").Value = 766.2
$ws.Range("J18").Value = 602
$ws.Range("L18").Value = 602
$ws.Range("N18").Value = -1170
# Row 70
$ws.Range("H70").Value = 140624.25
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
# Row 73
$ws.Range("H73").Value = 140624.25
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
# Row 112
$ws.Range("H112").Value = 1954.0769
$ws.Range("J112").Value = 2140.3
$ws.Range("L112").Value = 6420.900000000001
$ws.Range("N112").Value = -8636.900000000001
# Row 127
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()
# Row 129
$ws.Range("H129").Value = 5629
$ws.Range("I129").Value = 2898
$ws.Range("J129").Value = 6994.5
$ws.Range("K129").Value = 8694
$ws.Range("L129").Value = 20983.5
$ws.Range("M129").Value = -3694
$ws.Range("N129").Value = -30983.5
# Row 132
$ws.Range("H132").Value = 2988.6667
$ws.Range("I132").Value = 3016.9143
$ws.Range("K132").Value = 9050.742899999999
$ws.Range("M132").Value = -6520.742899999999
# Row 135
$ws.Range("H135").Value = 3596.2
$ws.Range("I135").Value = 2911.9167
$ws.Range("K135").Value = 26207.2503
$ws.Range("M135").Value = -23672.2503
# Row 137
$ws.Range("H137").Value = 2533.6206
$ws.Range("I137").Value = 1978.1305
$ws.Range("J137").Value = 4663
$ws.Range("K137").Value = 5934.3915
$ws.Range("L137").Value = 13989
$ws.Range("M137").Value = -3384.3915
$ws.Range("N137").Value = -19089
# Row 141
$ws.Range("H141").Value = 3476.5
$ws.Range("I141").Value = 3366.1714
$ws.Range("K141").Value = 10098.5142
$ws.Range("M141").Value = -4918.514200000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8967.3125
$ws.Range("I32").Value = 3575.3076
$ws.Range("K32").Value = 3575.3076
$ws.Range("M32").Value = -3288.3076
# Row 74
$ws.Range("H74").Value = 1843.9474
$ws.Range("I74").Value = 1846.2059
$ws.Range("K74").Value = 1846.2059
$ws.Range("M74").Value = -972.2058999999999
# Row 77
$ws.Range("H77").Value = 1843.9474
$ws.Range("I77").Value = 1846.2059
$ws.Range("K77").Value = 9231.029500000001
$ws.Range("M77").Value = -4863.029500000001

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 12882996
$ws.Range("J86").Value = 5471.1
$ws.Range("L86").Value = 5471.1
$ws.Range("N86").Value = -7717.1
# Row 89
$ws.Range("H89").Value = 12882996
$ws.Range("J89").Value = 5471.1
$ws.Range("L89").Value = 27355.5
$ws.Range("N89").Value = -38587.5
# Row 134
$ws.Range("H134").Value = 5005.0835
$ws.Range("I134").Value = 5005.0835
$ws.Range("K134").Value = 15015.2505
$ws.Range("M134").Value = -12480.2505

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("I31").Value = 4571.2915
$ws.Range("J31").Value = 2809.121
$ws.Range("K31").Value = 4571.2915
$ws.Range("L31").Value = 2809.121
$ws.Range("M31").Value = -4276.2915
$ws.Range("N31").Value = -3399.121
# Row 34
$ws.Range("I34").Value = 4571.2915
$ws.Range("J34").Value = 2809.121
$ws.Range("K34").Value = 4571.2915
$ws.Range("L34").Value = 2809.121
$ws.Range("M34").Value = -4369.2915
$ws.Range("N34").Value = -3213.121
# Row 107
$ws.Range("H107").Value = 753.625
$ws.Range("I107").Value = 571.63635
$ws.Range("J107").Value = 1154
$ws.Range("K107").Value = 571.63635
$ws.Range("L107").Value = 1154
$ws.Range("M107").Value = 1348.36365
$ws.Range("N107").Value = -4994
# Row 134
$ws.Range("H134").Value = 1521.1428
$ws.Range("I134").Value = 1484.3846
$ws.Range("K134").Value = 4453.1538
$ws.Range("M134").Value = -1918.1538
# Row 135
$ws.Range("H135").Value = 89949.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 89949.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 89949.5
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -100089.5

$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 499
$ws.Range("I9").Value = 499
$ws.Range("K9").Value = 499
$ws.Range("M9").Value = -329
# Row 12
$ws.Range("H12").Value = 100664.664
$ws.Range("I12").Value = 115797.6
$ws.Range("K12").Value = 115797.6
$ws.Range("M12").Value = -115657.6
# Row 17
$ws.Range("H17").Value = 222.5
$ws.Range("J17").Value = 245
$ws.Range("L17").Value = 245
$ws.Range("N17").Value = -581
# Row 96
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
# Row 102
$ws.Range("H102").Value = 2053.3333
$ws.Range("I102").Value = 2098.75
$ws.Range("K102").Value = 2098.75
$ws.Range("M102").Value = -476.75
# Row 113
$ws.Range("H113").Value = 995.3333
$ws.Range("I113").Value = 995.3333
$ws.Range("K113").Value = 995.3333
$ws.Range("M113").Value = 1174.6667
# Row 132
$ws.Range("H132").Value = 6047.067
$ws.Range("I132").Value = 6423.5
$ws.Range("J132").Value = 5294.2
$ws.Range("K132").Value = 19270.5
$ws.Range("L132").Value = 15882.6
$ws.Range("M132").Value = -16740.5
$ws.Range("N132").Value = -20942.6

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 35271.547
$ws.Range("I7").Value = 36998.7
$ws.Range("J7").Value = 18000
$ws.Range("K7").Value = 36998.7
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = -36886.7
$ws.Range("N7").Value = -18224
# Row 40
$ws.Range("H40").Value = 5332
$ws.Range("I40").Value = 4998.5
$ws.Range("J40").Value = 5999
$ws.Range("K40").Value = 4998.5
$ws.Range("L40").Value = 5999
$ws.Range("M40").Value = -4862.5
$ws.Range("N40").Value = -6271
# Row 126
$ws.Range("H126").Value = 35271.547
$ws.Range("I126").Value = 36998.7
$ws.Range("J126").Value = 18000
$ws.Range("K126").Value = 110996.1
$ws.Range("L126").Value = 54000
$ws.Range("M126").Value = -108526.1
$ws.Range("N126").Value = -58940

$ws = $wb.Worksheets.Item("WVR")
# Row 19
$ws.Range("H19").Value = 16253
$ws.Range("I19").Value = 2500
$ws.Range("K19").Value = 2500
$ws.Range("M19").Value = -2326
# Row 81
$ws.Range("H81").Value = 3575
$ws.Range("I81").Value = 6000
$ws.Range("K81").Value = 12000
$ws.Range("M81").Value = -10939
# Row 84
$ws.Range("H84").Value = 3575
$ws.Range("I84").Value = 6000
$ws.Range("K84").Value = 60000
$ws.Range("M84").Value = -54696
# Row 126
$ws.Range("H126").Value = 2499.5
$ws.Range("I126").Value = 2499.5
$ws.Range("K126").Value = 7498.5
$ws.Range("M126").Value = -5028.5
# Row 136
$ws.Range("H136").Value = 2679.1177
$ws.Range("I136").Value = 1632.0714
$ws.Range("K136").Value = 4896.2142
$ws.Range("M136").Value = -2346.2142
